$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")
$ws.Activate()

# --- Fill in previously-blank 0 values for weeks F/G/H on existing task rows ---
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0

$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0

$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0

$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0

$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0

$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0

$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0

$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0

$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0

$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0

$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0

# --- Row 19 ("Search Product by name") gains a new task label + hours ---
$ws.Range("D19").Value = "Testa"
$ws.Range("E19").Value = 6
$ws.Range("F19").Value = 6
$ws.Range("G19").Value = 6
$ws.Range("H19").Value = 6

# --- Insert 4 new task rows (new "floggit" DB user related tasks) above the totals ---
$ws.Range("A21:A24").EntireRow.Insert()

# copy the formatting used by the other data rows onto the freshly inserted rows
$ws.Range("A18:O18").Copy()
$ws.Range("A21:O24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A21").Value = "ADI floggit"
$ws.Range("B21").Value = "KFUAAA"
$ws.Range("C21").Value = "Task 6.1"
$ws.Range("D21").Value = "Git"
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = 9
$ws.Range("G21").Value = 9
$ws.Range("H21").Value = 9

$ws.Range("A22").Value = "ADI floggit"
$ws.Range("B22").Value = "KFUAAA"
$ws.Range("C22").Value = "Task 6.2"
$ws.Range("D22").Value = "Log4J"
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = 9

$ws.Range("A23").Value = "ADI floggit"
$ws.Range("B23").Value = "KFUAAA"
$ws.Range("C23").Value = "Task 6.3"
$ws.Range("D23").Value = "Junit"
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 9

$ws.Range("A24").Value = "ADI floggit"
$ws.Range("B24").Value = "KFUAAA"
$ws.Range("C24").Value = "Task 6.4"
$ws.Range("D24").Value = "Scrum Agile XP "
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 9
$ws.Range("G24").Value = 9
$ws.Range("H24").Value = 9

# --- Update the "Actual Burndown" totals row (now row 25) for the extended ranges ---
$ws.Range("E25").Formula = "=SUM(E3:E24)"
$ws.Range("F25").Formula = "=SUM(F3:F24)"
$ws.Range("G25").Formula = "=SUM(G3:G24)"
$ws.Range("H25").Formula = "=SUM(H10:H24)"

# --- Move the burndown chart down so it still sits below the (now longer) table ---
$co = $ws.ChartObjects(1)
$co.Top = $co.Top + 56

# --- Point the chart's "Actual burndown" series at the relocated totals row ---
$chart = $co.Chart
$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES(,,'Sprint 2'!`$E`$25:`$O`$25,1)"

# --- Restore the view: scrolled down a bit, with the last-edited cell selected ---
$ws.Range("I25").Select()
